$d = $word.ActiveDocument

# --- Starting layout ----------------------------------------------------
# Para 1: "Aprendiendo GIT el 11-02-18"
# Para 2: run "¿Hola cómo estás" + bookmark "_GoBack" + run "?"
#
# --- Target layout --------------------------------------------------------
# Para 1: unchanged
# Para 2: single run "¿Hola cómo estás?"
# Para 3 (new): single run "Todo bien, y vos?" followed by bookmark "_GoBack"
# ---------------------------------------------------------------------------

# Step 1: remove the "_GoBack" bookmark; it currently sits between
# "estás" and the trailing "?" run. It gets re-created later, after the
# newly-added sentence.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# Step 2: locate the sentence "¿Hola cómo estás" so we don't depend on
# hard-coded character offsets.
$findRng = $d.Content
$found = $findRng.Find.Execute("¿Hola cómo estás", $true, $false, $false, `
                                $false, $false, $true, 1, $false, "", 0)
$sentenceEnd = $findRng.End

# The trailing "?" lives in its own run right after that sentence.
# Delete that single character...
$qRng = $d.Range($sentenceEnd, $sentenceEnd + 1)
$qRng.Delete()

# ...and re-insert it directly onto the end of the sentence's run so Word
# merges it into "¿Hola cómo estás" -> "¿Hola cómo estás?" as one run.
$startPos = $sentenceEnd - 16
$r = $d.Range($startPos, $sentenceEnd)
$r.InsertAfter("?")

# Step 3: split the paragraph right after the new "¿Hola cómo estás?" so
# the following sentence becomes its own paragraph.
$splitPos = $sentenceEnd + 1
$splitRng = $d.Range($splitPos, $splitPos)
$splitRng.InsertParagraphAfter()

# Step 4: fill the (now empty) third paragraph with the new sentence.
$newParaStart = $splitPos + 1
$p3StartRng = $d.Range($newParaStart, $newParaStart)
$p3StartRng.InsertAfter("Todo bien, y vos?")

# Step 5: re-create the "_GoBack" bookmark at the end of the new
# paragraph, after its text. Adding a bookmark exactly at the end of the
# document's last paragraph is unreliable, so temporarily pad the story
# with one extra character, anchor the bookmark just before it, then
# remove the padding again.
$docEnd = $d.Content.End
$tailPos = $docEnd - 1
$tailRng = $d.Range($tailPos, $tailPos)
$tailRng.InsertAfter("X")
$bmRng = $d.Range($tailPos, $tailPos)
$d.Bookmarks.Add("_GoBack", $bmRng)
$padRng = $d.Range($tailPos, $tailPos + 1)
$padRng.Delete()
